$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.926.81'
$ws.Range("E2").Value = '  +2.00%  '

$ws.Range("D3").Value = '3.461.18'
$ws.Range("E3").Value = '  +0.79%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '3.462.61'
$ws.Range("E8").Value = '  +0.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.581'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.37'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.65%  '

$ws.Range("E11").Value = '  +2.74%  '

$ws.Range("E12").Value = '  +0.78%  '

$ws.Range("D13").Value = '4.058.34'
$ws.Range("E13").Value = '  +0.78%  '

$ws.Range("E14").Value = '  -2.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000195'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.33%  '

$ws.Range("D17").Value = '64.920.11'
$ws.Range("E17").Value = '  +1.78%  '

$ws.Range("D18").Value = '3.435.40'
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("E19").Value = '  -0.91%  '

$ws.Range("E20").Value = '  +0.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.64%  '

$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("E25").Value = '  +0.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000125'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +16.38%  '

$ws.Range("E27").Value = '  +1.38%  '

$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("E30").Value = '  +8.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.72'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.05%  '

$ws.Range("E34").Value = '  -1.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.10'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.60%  '

$ws.Range("E37").Value = '  +0.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.17'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.89%  '

$ws.Range("D40").Value = '3.012.91'
$ws.Range("E40").Value = '  +2.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0764'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.30'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.69%  '

$ws.Range("E43").Value = '  +4.81%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.70%  '

$ws.Range("E45").Value = '  -1.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.776'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.90%  '

$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.877'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.76%  '

$ws.Range("E50").Value = '  +3.38%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.62%  '
